$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.474.31'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.249.54'
$ws.Range("E3").Value = '  -2.88%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '497.45'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.48'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.296.95'
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0954'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("E11").Value = '  +2.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.324'
$ws.Range("E12").Value = '  +2.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.64'
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.668.26'
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.87'
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.455.69'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000130'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.306.42'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.06'
$ws.Range("E19").Value = '  +3.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.08'
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '305.85'
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.46'
$ws.Range("E22").Value = '  +3.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -1.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.99'
$ws.Range("E25").Value = '  -2.04%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.152'
$ws.Range("E27").Value = '  +6.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.374'
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.397.93'
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.16'
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.28'
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0691'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.93'
$ws.Range("E34").Value = '  +3.04%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.08'
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.991'
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.67'
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.21'
$ws.Range("E39").Value = '  +3.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.867'
$ws.Range("E40").Value = '  +1.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.68'
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.43'
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.41'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.376'
$ws.Range("E44").Value = '  +1.94%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '128.68'
$ws.Range("E46").Value = '  +3.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.91'
$ws.Range("E47").Value = '  +3.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0894'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.550'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '243.70'
$ws.Range("E50").Value = '  +1.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0484'
$ws.Range("E51").Value = '  +1.88%  '
